$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range('D2').Value = '42.937.05'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '2.336.15'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue $ws.Range('D5') '306.16'
$ws.Range('E5').Value = '  -1.75%  '
Set-TextValue $ws.Range('D6') '101.16'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').Value = '  -4.63%  '
$ws.Range('E8').Value = '  -0.01%  '
Set-TextValue $ws.Range('D9') '0.510'
$ws.Range('E9').Value = '  -3.30%  '
Set-TextValue $ws.Range('D10') '34.99'
$ws.Range('E10').Value = '  -2.32%  '
Set-TextValue $ws.Range('D11') '52.22'
$ws.Range('E11').Value = '  +1.33%  '
Set-TextValue $ws.Range('D12') '0.0799'
$ws.Range('E12').Value = '  -1.93%  '
Set-TextValue $ws.Range('D13') '0.113'
$ws.Range('E13').Value = '  +0.11%  '
Set-TextValue $ws.Range('D14') '6.81'
$ws.Range('E14').Value = '  -2.75%  '
Set-TextValue $ws.Range('D15') '15.91'
$ws.Range('E15').Value = '  +5.96%  '
$ws.Range('D16').Value = '2.328.75'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '42.865.54'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '0.0₃0909'
$ws.Range('E20').Value = '  -2.71%  '
Set-TextValue $ws.Range('D21') '11.67'
$ws.Range('E21').Value = '  -5.51%  '
Set-TextValue $ws.Range('D22') '67.90'
$ws.Range('E22').Value = '  -0.30%  '
Set-TextValue $ws.Range('D23') '236.78'
$ws.Range('E23').Value = '  -1.97%  '
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('E26').Value = '  -0.07%  '
Set-TextValue $ws.Range('D27') '25.39'
$ws.Range('E27').Value = '  +2.97%  '
$ws.Range('E28').Value = '  +9.64%  '
Set-TextValue $ws.Range('D29') '35.08'
$ws.Range('E29').Value = '  -5.48%  '
$ws.Range('E30').Value = '  -2.63%  '
Set-TextValue $ws.Range('D31') '160.93'
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  -3.05%  '
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('E35').Value = '  +6.10%  '
Set-TextValue $ws.Range('D36') '17.47'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('E38').Value = '  -4.51%  '
Set-TextValue $ws.Range('D39') '1.85'
$ws.Range('E39').Value = '  -1.75%  '
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('E41').Value = '  -2.39%  '
Set-TextValue $ws.Range('D42') '2.50'
$ws.Range('E42').Value = '  +7.75%  '
$ws.Range('D43').Value = '2.007.85'
$ws.Range('E43').Value = '  +1.69%  '
$ws.Range('E44').Value = '  -1.53%  '
Set-TextValue $ws.Range('D45') '18.72'
$ws.Range('E45').Value = '  -3.30%  '
Set-TextValue $ws.Range('D46') '10.19'
$ws.Range('E46').Value = '  +2.83%  '
Set-TextValue $ws.Range('D47') '2.93'
$ws.Range('E47').Value = '  -1.55%  '
Set-TextValue $ws.Range('D48') '55.68'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '2.562.32'
$ws.Range('E50').Value = '  +0.95%  '
$ws.Range('E51').Value = '  +3.28%  '
